$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-08-04 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-05 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("881÷2=440, 1", $true, $false, $false, $false, $false, $true, 1, $false, "101÷5=20, 1", 2) | Out-Null
$d.Content.Find.Execute("148÷5=29, 3", $true, $false, $false, $false, $false, $true, 1, $false, "855÷4=213, 3", 2) | Out-Null
$d.Content.Find.Execute("436÷5=87, 1", $true, $false, $false, $false, $false, $true, 1, $false, "935÷5=187, 0", 2) | Out-Null
$d.Content.Find.Execute("640÷4=160, 0", $true, $false, $false, $false, $false, $true, 1, $false, "822÷8=102, 6", 2) | Out-Null
$d.Content.Find.Execute("918÷4=229, 2", $true, $false, $false, $false, $false, $true, 1, $false, "585÷4=146, 1", 2) | Out-Null
$d.Content.Find.Execute("652÷4=163, 0", $true, $false, $false, $false, $false, $true, 1, $false, "410÷4=102, 2", 2) | Out-Null
$d.Content.Find.Execute("915÷2=457, 1", $true, $false, $false, $false, $false, $true, 1, $false, "603÷6=100, 3", 2) | Out-Null
$d.Content.Find.Execute("494÷3=164, 2", $true, $false, $false, $false, $false, $true, 1, $false, "271÷7=38, 5", 2) | Out-Null
$d.Content.Find.Execute("320÷3=106, 2", $true, $false, $false, $false, $false, $true, 1, $false, "885÷8=110, 5", 2) | Out-Null
$d.Content.Find.Execute("212÷7=30, 2", $true, $false, $false, $false, $false, $true, 1, $false, "496÷5=99, 1", 2) | Out-Null
$d.Content.Find.Execute("901÷5=180, 1", $true, $false, $false, $false, $false, $true, 1, $false, "128÷2=64, 0", 2) | Out-Null
$d.Content.Find.Execute("672÷8=84, 0", $true, $false, $false, $false, $false, $true, 1, $false, "799÷9=88, 7", 2) | Out-Null
$d.Content.Find.Execute("361÷9=40, 1", $true, $false, $false, $false, $false, $true, 1, $false, "690÷3=230, 0", 2) | Out-Null
$d.Content.Find.Execute("538÷2=269, 0", $true, $false, $false, $false, $false, $true, 1, $false, "801÷7=114, 3", 2) | Out-Null
$d.Content.Find.Execute("700÷6=116, 4", $true, $false, $false, $false, $false, $true, 1, $false, "451÷2=225, 1", 2) | Out-Null
$d.Content.Find.Execute("302÷4=75, 2", $true, $false, $false, $false, $false, $true, 1, $false, "120÷3=40, 0", 2) | Out-Null
$d.Content.Find.Execute("289÷5=57, 4", $true, $false, $false, $false, $false, $true, 1, $false, "497÷4=124, 1", 2) | Out-Null
$d.Content.Find.Execute("487÷7=69, 4", $true, $false, $false, $false, $false, $true, 1, $false, "944÷9=104, 8", 2) | Out-Null
$d.Content.Find.Execute("201÷8=25, 1", $true, $false, $false, $false, $false, $true, 1, $false, "233÷9=25, 8", 2) | Out-Null
$d.Content.Find.Execute("978÷9=108, 6", $true, $false, $false, $false, $false, $true, 1, $false, "248÷8=31, 0", 2) | Out-Null
$d.Content.Find.Execute("379÷9=42, 1", $true, $false, $false, $false, $false, $true, 1, $false, "321÷4=80, 1", 2) | Out-Null
$d.Content.Find.Execute("388÷5=77, 3", $true, $false, $false, $false, $false, $true, 1, $false, "591÷7=84, 3", 2) | Out-Null
$d.Content.Find.Execute("733÷9=81, 4", $true, $false, $false, $false, $false, $true, 1, $false, "766÷5=153, 1", 2) | Out-Null
$d.Content.Find.Execute("771÷7=110, 1", $true, $false, $false, $false, $false, $true, 1, $false, "160÷2=80, 0", 2) | Out-Null
$d.Content.Find.Execute("717÷7=102, 3", $true, $false, $false, $false, $false, $true, 1, $false, "390÷2=195, 0", 2) | Out-Null
